# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet at
# position N (14) - this is the "Late" grouping gap column that shifts the
# existing Late/heading/Outstanding columns one place to the right - and
# make "Repayment schedule" the active/selected sheet & cell, matching the
# workbook being reopened with that tab focused.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14), shifting N,O,P -> O,P,Q
$ws.Columns.Item(14).Insert()

# New column gets width 11 (no bestFit, just a custom width) like the rest
# of the block next to it.
$ws.Columns.Item(14).ColumnWidth = 10.2

# Move the selection / active cell to R6 and make this sheet the active tab
# (this also clears tabSelected on whatever sheet was active before).
$ws.Activate()
$ws.Range("R6").Select()
